# Add four new days (Mon 4/15 - Thu 4/18 2013) to the bottom of the workload
# log on Sheet1, finishing with a merged "away" note row for the Brisbane
# studio trip (mirrors the existing merged rows, e.g. B56:E56 "Away").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 72: Mon 2013-04-15 (serial 41379) ---
[void]$ws.Range("A71:E71").Copy()
[void]$ws.Range("A72:E72").PasteSpecial(-4122)
$ws.Range("A72").Value = 41379
$ws.Range("B72").Value = "1h"
$ws.Range("C72").Value = "0H"
$ws.Range("D72").Value = "1H"
$ws.Range("E72").Value = "-"

# --- Row 73: Tue 2013-04-16 (serial 41380) ---
[void]$ws.Range("A71:E71").Copy()
[void]$ws.Range("A73:E73").PasteSpecial(-4122)
$ws.Range("A73").Value = 41380
$ws.Range("B73").Value = "0H"
$ws.Range("C73").Value = "0H"
$ws.Range("D73").Value = "1H"
$ws.Range("E73").Value = "-"

# --- Row 74: Wed 2013-04-17 (serial 41381) ---
[void]$ws.Range("A71:E71").Copy()
[void]$ws.Range("A74:E74").PasteSpecial(-4122)
$ws.Range("A74").Value = 41381
$ws.Range("B74").Value = "0H"
$ws.Range("C74").Value = "0H"
$ws.Range("D74").Value = "1H"
$ws.Range("E74").Value = "-"

# --- Row 75: Thu 2013-04-18 (serial 41382) - merged "away" note row ---
[void]$ws.Range("A56:E56").Copy()
[void]$ws.Range("A75:E75").PasteSpecial(-4122)
$ws.Range("A75").Value = 41382
$ws.Range("B75").Value = "BRISBANE TRIP TO THE STUDIOS"
[void]$ws.Range("B75:E75").Merge()

$excel.CutCopyMode = $false

# --- Update the view: scroll position + active selection ---
$ws.Activate()
[void]$ws.Range("G53").Select()
$excel.ActiveWindow.ScrollRow = 37
